# Completion log update: mark additional upper-case letters as "Completed".
#
# Rows 32-57 hold letters A-Z in column B, with a Yes/No "Completed" flag in
# column C (green "Yes" style / red "No" style) and an optional completion
# date in column D.
#
# This change flips all of A-Z to "Yes" except L, O and X (rows 43, 46, 55),
# which are left as "No" per the commit message ("will need to supplement a
# few characters"). Rows H and I (39, 40) additionally get a completion date
# of 2019-06-11 (serial 43627).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A "Yes" cell (green fill) to copy formatting from, and a "date" cell to
# copy the date-formatted style from.
$yesFormatSource = $ws.Range("C6")
$dateFormatSource = $ws.Range("D6")

$rowsToSkip = @(43, 46, 55)
$rowsWithDate = @(39, 40)

for ($row = 32; $row -le 57; $row++) {
    if ($rowsToSkip -contains $row) {
        continue
    }

    $cCell = $ws.Cells.Item($row, 3)
    $yesFormatSource.Copy()
    $cCell.PasteSpecial(-4122)  # xlPasteFormats
    $cCell.Value = "Yes"

    if ($rowsWithDate -contains $row) {
        $dCell = $ws.Cells.Item($row, 4)
        $dateFormatSource.Copy()
        $dCell.PasteSpecial(-4122)  # xlPasteFormats
        $dCell.Value = 43627
    }
}

$excel.CutCopyMode = $false

# Update the saved view state (scroll position + active selection) to match
# where the author was working when they made this edit.
$activeWindow = $excel.ActiveWindow
$activeWindow.ScrollRow = 33
$activeWindow.ScrollColumn = 1
$ws.Range("D60").Select()
